$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 75, shifting existing rows 75-201 down to 76-202
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new data
$ws.Cells.Item(75, 1).Value = 8
$ws.Cells.Item(75, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 44477
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = 100114013
$ws.Cells.Item(75, 7).Value = "Zanahoria"
$ws.Cells.Item(75, 8).Value = "Sin especificar"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 800
$ws.Cells.Item(75, 11).Value = 6000
$ws.Cells.Item(75, 12).Value = 7000
$ws.Cells.Item(75, 13).Value = 6500
$ws.Cells.Item(75, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(75, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(75, 16).Value = 325
$ws.Cells.Item(75, 17).Value = 20
$ws.Cells.Item(75, 18).Value = "Hortaliza"
